$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up the option_values text (shorter separators, Vietnamese wording tweak)
$ws.Range("E2").Value = "i5-6500T; Không RAM, Không SSD"
$ws.Range("E3").Value = "i5-6500T; Không RAM, 256GB NVMe 95%"
$ws.Range("E4").Value = "i5-6500T; Không RAM, 256GB NVMe"
$ws.Range("E5").Value = "i5-6500;, Không RAM, 512GB NVMe"

# Apply an explicit font to the option_values column (creates a second font/style entry)
$ws.Range("E2:E5").Font.ThemeFont = 1

# Match the final selection left behind in the sheet
$ws.Range("F11").Select()
